$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Auto-generated cell value updates (computed from diff)
$updates = @{
    "J25" = -8.180305474188032
    "K25" = -0.04062844669037713
    "I26" = -8.180305474188014
    "J26" = -0.04062844669036003
    "H27" = -8.130305474188006
    "I27" = 0.009371553309648506
    "G28" = -8.180305474188014
    "H28" = -0.04062844669036003
    "F29" = -8.131305474188011
    "G29" = 0.008371553309643787
    "H29" = 2.155942556679634
    "I29" = -1.284091849519598
    "J29" = -1.345275412033999
    "K29" = 0.6961888372093057
    "E30" = -8.150305474188002
    "F30" = -0.0006284466903566099
    "G30" = 2.206942556679635
    "H30" = -1.233091849519596
    "I30" = -1.314275412034007
    "J30" = 0.6871888372093053
    "D31" = -6.480305474188057
    "E31" = 1.359371553309629
    "F31" = 3.306942556679616
    "G31" = -0.6330918495195741
    "H31" = -1.312101506057118
    "I31" = 0.6891265057146824
    "C32" = -1.580305474188014
    "D32" = 3.45937155330964
    "E32" = 4.00694255667963
    "F32" = -0.533091849519602
    "G32" = -1.194275412034003
    "H32" = 0.7471888372093018
    "B33" = -1.840305474188014
    "C33" = 2.411492840269247
    "D33" = 3.1636778921727
    "E33" = -0.8080918495195992
    "F33" = -1.194275412034017
    "G33" = 0.7396061149258059
    "H33" = 0.1538456963687809
    "I33" = 0.2501940957558864
    "J33" = -0.009553352404694293
    "K33" = 0.8265826272929218
    "B34" = -0.7406284466903656
    "C34" = 2.395386534006718
    "D34" = 0.9669081504804069
    "E34" = -0.9882337193940176
    "F34" = 0.946940841384361
    "G34" = 0.228530305674387
    "H34" = 0.2998830883735715
    "I34" = 0.03126059032258927
    "J34" = 0.859248693144868
    "B35" = 1.097462556679631
    "C35" = 0.4233181504804118
    "D35" = 0.1547245879659871
    "E35" = 1.2036488372093
    "F35" = 0.2516179039985633
    "G35" = 0.3175051528550954
    "H35" = 0.0327531447941104
    "I35" = 0.8612811306276531
    "B36" = -0.323091849519586
    "C36" = 1.865724587965999
    "D36" = 1.847581740563939
    "E36" = 0.2513302348380089
    "F36" = 0.321830387954793
    "G36" = 0.03045639295955022
    "H36" = 0.8604233730577899
    "B37" = -0.1742754120339929
    "C37" = 1.647188837209285
    "D37" = 1.424597903998563
    "E37" = 1.030275152855095
    "F37" = 0.3098531447941104
    "G37" = 0.9253111306276531
    "H37" = 0.595288361530445
    "I37" = 0.2314303014965162
    "J37" = 0.3396233709318813
    "K37" = 0.4970570324401455
    "B38" = -0.1928111627906902
    "C38" = -0.2454020960014367
    "D38" = 2.200275152855095
    "E38" = 1.34985314479411
    "F38" = 1.245311130627653
    "G38" = 0.865288361530445
    "H38" = 0.3014303014965162
    "I38" = 0.3896233709318813
    "J38" = 0.5470570324401456
    "B39" = 0.2298834519089068
    "C39" = 0.852673321422694
    "D39" = 0.6965172908002875
    "E39" = 1.5451969683358
    "F39" = 1.226054924247761
    "G39" = 0.7182574749033117
    "H39" = 0.6131599809998507
    "I39" = 0.6331062931918918
    "B40" = 0.2402751528550954
    "C40" = 0.7703367535273524
    "D40" = 1.835311130627653
    "E40" = 1.375288361530445
    "F40" = 0.8211776578870573
    "G40" = 0.7472245760093372
    "H40" = 0.7070570324401455
    "B41" = -0.3601468552058896
    "C41" = 0.2553111306276531
    "D41" = -0.06471163846955513
    "E41" = -0.01856969850348378
    "F41" = 0.6096233709318812
    "G41" = 0.7670570324401456
    "H41" = 0.2964009461638233
    "I41" = 0.5467572334344299
    "J41" = 0.3347911218750448
    "K41" = 0.6205382402049349
    "B42" = 0.1453111306276531
    "C42" = -0.05471163846955512
    "D42" = 0.05096192991250348
    "E42" = 0.6813638604798342
    "F42" = 0.7891405015358686
    "G42" = 0.2056981497352126
    "H42" = 0.4324739641823357
    "I42" = 0.1941639136379365
    "J42" = 0.4450409605369429
    "B43" = 0.1274104041151531
    "C43" = 0.04629456727128911
    "D43" = 0.5849680602752159
    "E43" = 0.8484308539411956
    "F43" = 0.3127969308147129
    "G43" = 0.4924679903592591
    "H43" = 0.1829809721025981
    "I43" = 0.4311701029554768
    "B44" = 0.1020781346750965
    "C44" = 0.3096233709318813
    "D44" = 0.7070570324401455
    "E44" = 0.2864009461638233
    "F44" = 0.5267572334344299
    "G44" = 0.2247911218750448
    "H44" = 0.5305382402049349
    "B45" = -0.05037662906811868
    "C45" = 0.6070570324401456
    "D45" = 0.2664009461638233
    "E45" = 0.6279049279394684
    "F45" = 0.3979720162165136
    "G45" = 0.5812217942731905
    "H45" = -0.07988875890465857
    "I45" = 0.5355197892996415
    "B46" = 0.3310570324401456
    "C46" = 0.08037229097246271
    "D46" = 0.4669936816748645
    "E46" = 0.3339657216846063
    "F46" = 0.6406675981017713
    "G46" = -0.0121054429918957
    "H46" = 0.5766843069119603
    "B47" = -0.2961090942996805
    "C47" = 0.2811150946827183
    "D47" = 0.4061476709679112
    "E47" = 0.7325305045040693
    "F47" = 0.03215161503959513
    "G47" = 0.6272170465139766
    "B48" = 0.3250097553919601
    "C48" = 0.3247911218750448
    "D48" = 0.720538240204935
    "E48" = 0.04834807457247964
    "F48" = 0.6666843069119602
    "B49" = -0.1523844033498989
    "C49" = 0.4105382402049349
    "D49" = -0.09865192542752038
    "E49" = 0.6966843069119601
    "B50" = 0.4147680207538826
    "C50" = -0.211887816642232
    "D50" = 0.4271580360664302
    "B51" = -0.1722147680129069
    "C51" = 0.4775246863929397
    "B52" = 0.3238631410950035
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Cells that were removed/cleared entirely in the edit
$clears = @(
    "J45",
    "I46",
    "H47",
    "G48",
    "F49",
    "E50",
    "D51",
    "C52",
    "B53",
)
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}
